# Decrement the "剩余" (remaining days) column E for each data row by 1,
# mirroring the daily countdown update captured in the diff.
# Row 36 (E36) is intentionally left untouched, matching the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 5).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 99 }

for ($r = 2; $r -le $lastRow; $r++) {
    if ($r -eq 36) { continue }
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value2 = $val - 1
    }
}
